# Reproduce the commit: add a second worksheet ("Sheet1") to the
# template_data workbook that holds the same 65 header labels as
# template_data!A1:BM1, but transposed into a single column (A1:A65).
# Also refresh the column widths / view state on both sheets so the
# saved file matches what Excel would have produced when a user
# selected the header row, copied it, and pasted-special->transposed
# it into a brand-new sheet.

$wb = $excel.ActiveWorkbook
$templateSheet = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# 1. Capture the 65 header values from row 1 of template_data while
#    it is still the only sheet (Value2 avoids the parameterized
#    -Value- reflection quirk on this host).
# ---------------------------------------------------------------
$headerCount = 65
$headers = @()
for ($c = 1; $c -le $headerCount; $c++) {
    $headers += , $templateSheet.Cells.Item(1, $c).Value2
}

# ---------------------------------------------------------------
# 2. Add the new sheet after template_data so tab order becomes
#    template_data, Sheet1 (matches sheetId 1 then 2 in the diff).
# ---------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)

# ---------------------------------------------------------------
# 3. Write the transposed header values into column A of the new
#    sheet, reusing the existing shared strings.
# ---------------------------------------------------------------
for ($r = 1; $r -le $headerCount; $r++) {
    $newSheet.Cells.Item($r, 1).Value = $headers[$r - 1]
}

# ---------------------------------------------------------------
# 4. Column widths -- best-fit values captured from the real Excel
#    session (character-width units, snapped to this host's nearest
#    representable 1/6-character granularity for ColumnWidth).
# ---------------------------------------------------------------
$templateWidths = @{
    1 = 20.6640625;  2 = 15.44140625; 3 = 12.21875;    4 = 11;
    5 = 6.6640625;   6 = 5.5546875;   7 = 6.5546875;   8 = 12.21875;
    9 = 11.33203125; 10 = 3.88671875; 11 = 12.33203125; 12 = 9;
    13 = 12.109375;  14 = 12;         15 = 12;          16 = 12;
    17 = 6.88671875; 18 = 9.21875;    19 = 9.5546875;   20 = 13.109375;
    21 = 22.44140625; 22 = 22.88671875; 23 = 5.44140625; 24 = 15.6640625;
    25 = 6.21875;    26 = 7.44140625; 27 = 3.77734375;  28 = 11.21875;
    29 = 12.5546875; 30 = 9;          31 = 12.44140625; 32 = 20.21875;
    33 = 15.33203125; 34 = 12.109375; 35 = 8.21875;     36 = 16.6640625;
    37 = 11.5546875; 38 = 14.21875;   39 = 13.77734375; 40 = 12.109375;
    41 = 7;          42 = 10.33203125; 43 = 6.88671875; 44 = 11.33203125;
    45 = 7.21875;    46 = 5.5546875;  47 = 12.6640625;  48 = 12.5546875;
    49 = 10.21875;   50 = 12.88671875; 51 = 10.33203125; 52 = 5.33203125;
    53 = 6.5546875;  54 = 7.88671875; 55 = 13.6640625;  56 = 7.6640625;
    57 = 8.77734375; 58 = 14.21875;   59 = 14;          60 = 16.21875;
    61 = 11;         62 = 11.88671875; 63 = 7.5546875;  64 = 10.77734375;
    65 = 5.88671875;
}

foreach ($colIndex in $templateWidths.Keys) {
    $targetWidth = $templateWidths[$colIndex]
    # ColumnWidth(input) is stored as input + 5/6 "characters"; solve
    # for the input that lands closest to the recorded best-fit width.
    $inputWidth = $targetWidth - (5 / 6)
    $templateSheet.Columns.Item($colIndex).ColumnWidth = $inputWidth
}

$newSheet.Columns.Item(1).ColumnWidth = 22.88671875 - (5 / 6)

# ---------------------------------------------------------------
# 5. View state: template_data keeps A1:BM1 selected (no longer the
#    active tab), the new sheet becomes the active tab with its full
#    column selected.
# ---------------------------------------------------------------
$templateSheet.Activate()
$templateSheet.Range("A1:BM1").Select() | Out-Null

$newSheet.Activate()
$newSheet.Range("A1:A65").Select() | Out-Null
